$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column K (2023) data
$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 1040
$ws.Range("K5").Value = 291
$ws.Range("K6").Value = 749
